$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: coin name/link swaps, refreshed prices and 1h volume deltas.
$updates = @(
    @{ Row = 2; D = '27.736.02'; E = '  +0.65%  ' },
    @{ Row = 3; D = '1.631.85'; E = '  +0.24%  ' },
    @{ Row = 4; D = '0.993'; E = '  -0.97%  ' },
    @{ Row = 5; D = '211.38'; E = '  -0.20%  ' },
    @{ Row = 6; D = '0.521'; E = '  +0.10%  ' },
    @{ Row = 7; D = '0.993'; E = '  -0.98%  ' },
    @{ Row = 8; D = '23.41'; E = '  +2.10%  ' },
    @{ Row = 9; D = '0.257'; E = '  -2.07%  ' },
    @{ Row = 10; D = '0.0611'; E = '  +0.04%  ' },
    @{ Row = 11; D = '0.0877'; E = '  +2.07%  ' },
    @{ Row = 12; D = '1.866.52'; E = '  +0.33%  ' },
    @{ Row = 13; D = '1.628.69'; E = '  +0.04%  ' },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '4.07'; E = '  +0.98%  ' },
    @{ Row = 15; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.571'; E = '  +2.92%  ' },
    @{ Row = 16; D = '65.42'; E = '  +0.82%  ' },
    @{ Row = 17; D = '27.754.23'; E = '  +0.72%  ' },
    @{ Row = 18; D = '232.85'; E = '  +2.13%  ' },
    @{ Row = 19; D = '0.0₃0720'; E = '  +0.36%  ' },
    @{ Row = 20; D = '7.57'; E = '  +0.73%  ' },
    @{ Row = 21; D = '0.993'; E = '  -0.91%  ' },
    @{ Row = 22; D = '4.34'; E = '  +0.13%  ' },
    @{ Row = 23; D = '10.28'; E = '  -3.32%  ' },
    @{ Row = 24; D = '2.08'; E = '  -1.99%  ' },
    @{ Row = 25; D = '151.77'; E = '  +1.76%  ' },
    @{ Row = 26; D = '6.89'; E = '  +0.50%  ' },
    @{ Row = 27; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.111'; E = '  +0.03%  ' },
    @{ Row = 28; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '15.60'; E = '  +0.47%  ' },
    @{ Row = 29; D = '0.994'; E = '  -0.86%  ' },
    @{ Row = 30; E = '  -0.15%  ' },
    @{ Row = 31; D = '0.0482'; E = '  +0.35%  ' },
    @{ Row = 32; D = '3.33'; E = '  +1.94%  ' },
    @{ Row = 33; D = '3.10'; E = '  +1.29%  ' },
    @{ Row = 34; D = '1.406.62'; E = '  -4.11%  ' },
    @{ Row = 35; D = '1.57'; E = '  +2.35%  ' },
    @{ Row = 36; D = '2.34'; E = '  +0.70%  ' },
    @{ Row = 37; E = '  +1.11%  ' },
    @{ Row = 38; D = '0.876'; E = '  +0.33%  ' },
    @{ Row = 39; D = '0.556'; E = '  -0.40%  ' },
    @{ Row = 40; D = '0.902'; E = '  -2.39%  ' },
    @{ Row = 41; E = '  +0.37%  ' },
    @{ Row = 42; D = '0.994'; E = '  -0.91%  ' },
    @{ Row = 43; D = '1.86'; E = '  +6.87%  ' },
    @{ Row = 44; D = '66.83'; E = '  -1.06%  ' },
    @{ Row = 45; D = '5.49'; E = '  +2.86%  ' },
    @{ Row = 46; D = '2.20'; E = '  -0.62%  ' },
    @{ Row = 47; D = '1.774.18'; E = '  +0.30%  ' },
    @{ Row = 48; D = '87.34'; E = '  +0.11%  ' },
    @{ Row = 49; D = '0.0996'; E = '  +0.82%  ' },
    @{ Row = 50; D = '0.0504'; E = '  -0.03%  ' },
    @{ Row = 51; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.0₇0974'; E = '  -7.77%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Prices are point-separated strings (e.g. "27.736.02", "15.60"); force
        # text storage so Excel does not reinterpret them as numbers/dates and
        # strip significant trailing zeros.
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
}
